# Edit the pledged-item reference number on Sheet1!A2.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Select the cell being edited (mirrors the natural click-then-type flow
# that moved the saved cursor position from B13 to A2).
$ws.Range("A2").Select()

# Update the cell's text value (the actual content change of this commit).
$ws.Range("A2").Value = "98611120170160211-01333"
